# The "Statistics" sheet lists one row per university profile (Profile,
# AVG Score, Number of Universities, Number of Students, Universities names).
# This edit reorders the data rows so the MEDICINE profile moves from the
# bottom (row 5) up to the top of the data (row 2), shifting LINGUISTICS and
# MATHEMATICS down by one row each. PHYSICS (row 3) is unaffected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: now holds the MEDICINE profile (previously on row 5)
$ws.Cells.Item(2, 1).Value2 = "MEDICINE"
$ws.Cells.Item(2, 2).Value2 = 4.329999923706055
$ws.Cells.Item(2, 3).Value2 = 3.0
$ws.Cells.Item(2, 4).Value2 = 3.0
$ws.Cells.Item(2, 5).Value2 = "МГМУ;ТУМ;СМИ"

# Row 3: PHYSICS is unchanged

# Row 4: now holds the LINGUISTICS profile (previously on row 2)
$ws.Cells.Item(4, 1).Value2 = "LINGUISTICS"
$ws.Cells.Item(4, 2).Value2 = 0.0
$ws.Cells.Item(4, 3).Value2 = 1.0
$ws.Cells.Item(4, 4).Value2 = 0.0
$ws.Cells.Item(4, 5).Value2 = "ВЛПУ"

# Row 5: now holds the MATHEMATICS profile (previously on row 4)
$ws.Cells.Item(5, 1).Value2 = "MATHEMATICS"
$ws.Cells.Item(5, 2).Value2 = 0.0
$ws.Cells.Item(5, 3).Value2 = 1.0
$ws.Cells.Item(5, 4).Value2 = 0.0
$ws.Cells.Item(5, 5).Value2 = "КУВ"
